# New PO forecast model
# Updates three sheets:
#  - "Weekly Quantity": append 3 new weekly rows (33-35)
#  - "Monthly Trend": append 1 new monthly row (10)
#  - "PO Forecast": refresh the whole forecast series (new model values for
#    existing rows 2-32) and extend the series with new forecast rows
#    (33-43), shifting the remaining weeks forward.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Weekly Quantity  (dimension A1:B32 -> A1:B35)
# ---------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$weeklyDateFmt = $wsWeekly.Range("A32").NumberFormat

$weeklyNewRows = @(
    @{ Row = 33; Date = 45662.99999999999; Qty = 74 },
    @{ Row = 34; Date = 45669.99999999999; Qty = 26 },
    @{ Row = 35; Date = 45676.99999999999; Qty = 2 }
)

foreach ($r in $weeklyNewRows) {
    $wsWeekly.Cells.Item($r.Row, 1).Value = $r.Date
    $wsWeekly.Cells.Item($r.Row, 1).NumberFormat = $weeklyDateFmt
    $wsWeekly.Cells.Item($r.Row, 2).Value = $r.Qty
}

# ---------------------------------------------------------------------
# Sheet 2: Monthly Trend  (dimension A1:B9 -> A1:B10)
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$monthlyDateFmt = $wsMonthly.Range("A9").NumberFormat

$wsMonthly.Cells.Item(10, 1).Value = 45688.99999999999
$wsMonthly.Cells.Item(10, 1).NumberFormat = $monthlyDateFmt
$wsMonthly.Cells.Item(10, 2).Value = 102

# ---------------------------------------------------------------------
# Sheet 3: PO Forecast  (dimension A1:B40 -> A1:B43)
# ---------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")
$forecastDateFmt = $wsForecast.Range("A32").NumberFormat

$forecastRows = @(
    @{ Row = 2; Date = 45389.99999999999; Qty = 64 },
    @{ Row = 3; Date = 45396.99999999999; Qty = 64 },
    @{ Row = 4; Date = 45403.99999999999; Qty = 64 },
    @{ Row = 5; Date = 45410.99999999999; Qty = 63 },
    @{ Row = 6; Date = 45417.99999999999; Qty = 63 },
    @{ Row = 7; Date = 45424.99999999999; Qty = 63 },
    @{ Row = 8; Date = 45431.99999999999; Qty = 63 },
    @{ Row = 9; Date = 45438.99999999999; Qty = 63 },
    @{ Row = 10; Date = 45445.99999999999; Qty = 63 },
    @{ Row = 11; Date = 45452.99999999999; Qty = 62 },
    @{ Row = 12; Date = 45459.99999999999; Qty = 62 },
    @{ Row = 13; Date = 45466.99999999999; Qty = 62 },
    @{ Row = 14; Date = 45473.99999999999; Qty = 62 },
    @{ Row = 15; Date = 45480.99999999999; Qty = 62 },
    @{ Row = 16; Date = 45487.99999999999; Qty = 61 },
    @{ Row = 17; Date = 45494.99999999999; Qty = 61 },
    @{ Row = 18; Date = 45501.99999999999; Qty = 61 },
    @{ Row = 19; Date = 45508.99999999999; Qty = 61 },
    @{ Row = 20; Date = 45515.99999999999; Qty = 61 },
    @{ Row = 21; Date = 45522.99999999999; Qty = 61 },
    @{ Row = 22; Date = 45529.99999999999; Qty = 60 },
    @{ Row = 23; Date = 45536.99999999999; Qty = 60 },
    @{ Row = 24; Date = 45571.99999999999; Qty = 59 },
    @{ Row = 25; Date = 45578.99999999999; Qty = 59 },
    @{ Row = 26; Date = 45585.99999999999; Qty = 59 },
    @{ Row = 27; Date = 45606.99999999999; Qty = 59 },
    @{ Row = 28; Date = 45613.99999999999; Qty = 58 },
    @{ Row = 29; Date = 45620.99999999999; Qty = 58 },
    @{ Row = 30; Date = 45627.99999999999; Qty = 58 },
    @{ Row = 31; Date = 45634.99999999999; Qty = 58 },
    @{ Row = 32; Date = 45641.99999999999; Qty = 58 },
    @{ Row = 33; Date = 45662.99999999999; Qty = 57 },
    @{ Row = 34; Date = 45669.99999999999; Qty = 57 },
    @{ Row = 35; Date = 45676.99999999999; Qty = 57 },
    @{ Row = 36; Date = 45683.99999999999; Qty = 57 },
    @{ Row = 37; Date = 45690.99999999999; Qty = 56 },
    @{ Row = 38; Date = 45697.99999999999; Qty = 56 },
    @{ Row = 39; Date = 45704.99999999999; Qty = 56 },
    @{ Row = 40; Date = 45711.99999999999; Qty = 56 },
    @{ Row = 41; Date = 45718.99999999999; Qty = 56 },
    @{ Row = 42; Date = 45725.99999999999; Qty = 56 },
    @{ Row = 43; Date = 45732.99999999999; Qty = 55 }
)

foreach ($r in $forecastRows) {
    $wsForecast.Cells.Item($r.Row, 1).Value = $r.Date
    $wsForecast.Cells.Item($r.Row, 1).NumberFormat = $forecastDateFmt
    $wsForecast.Cells.Item($r.Row, 2).Value = $r.Qty
}

Write-Output "PO forecast model updated"
